$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ftests")

# Fill in row 23 with the new "rule 21" test case
$ws.Range("C23").Value = "% TIV deductible with min and max deductible"
$ws.Range("D23").Value = "All"
$ws.Range("E23").Value = 21
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = "in progress"
$ws.Range("I23").Value = "in progress"

# Update the active selection as recorded in the sheet view
$ws.Range("C26").Select()
